# Lean UX canvas - update the Creative Commons license reference in the
# speaker notes of slide 1 and slide 2 from "BY-NC-SA 3.0" to "BY-NC-SA 4.0"
# (with the Dutch deed page), as described by the commit:
#   "Link naar licentie aangepast obv andere NL producten (v4 van CC licentie)"
#
# The notes text body has to be rewritten as a whole (this host's TextRange
# only supports whole-frame .Text assignment), so we read the current text,
# patch only the license sentence, and write the full text back together -
# rebuilding paragraph breaks with `n so the paragraph count/order is kept.

$p = $ppt.ActivePresentation

function Update-LicenseNotes {
    param($slideIndex)

    $slide = $p.Slides.Item($slideIndex)
    $notesShape = $slide.NotesPage.Shapes.Item(2)
    $tr = $notesShape.TextFrame.TextRange

    $paragraphs = $tr.Text -split "`r"

    for ($i = 0; $i -lt $paragraphs.Count; $i++) {
        if ($paragraphs[$i] -like "*creativecommons.org*") {
            $paragraphs[$i] = $paragraphs[$i].Replace("/3.0/ ", "/4.0/deed.nl")
        }
    }

    $tr.Text = [string]::Join("`n", $paragraphs)
}

Update-LicenseNotes 1
Update-LicenseNotes 2
